$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The race name ("Epreuve") and its date changed for the event listed in rows 2-9,
# and the id_epreuve code moved from 2 to 4.
$range = $ws.Range("E2:F9")
$range.NumberFormat = "@"

for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 4).Value = 4
    $ws.Cells.Item($r, 5).Value = "soleil"
    $ws.Cells.Item($r, 6).Value = "2021-03-02"
}

$range.ClearFormats()
